$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.231.92'
$ws.Range('E2').Value = '  -5.06%  '

$ws.Range('D3').Value = '3.036.35'
$ws.Range('E3').Value = '  -7.78%  '

$ws.Range('E4').Value = '  +0.32%  '

$ws.Range('D5').Value = '554.52'
$ws.Range('E5').Value = '  -7.13%  '

$ws.Range('D6').Value = '141.36'
$ws.Range('E6').Value = '  -7.61%  '

$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('D8').Value = '3.028.09'
$ws.Range('E8').Value = '  -7.91%  '

$ws.Range('D9').Value = '0.487'
$ws.Range('E9').Value = '  -10.56%  '

$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').Value = '6.42'
$ws.Range('E10').Value = '  -4.71%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.154'
$ws.Range('E11').Value = '  -11.40%  '

$ws.Range('D12').Value = '0.464'
$ws.Range('E12').Value = '  -9.95%  '

$ws.Range('D13').Value = '36.46'
$ws.Range('E13').Value = '  -6.23%  '

$ws.Range('D14').Value = '0.0000221'
$ws.Range('E14').Value = '  -10.43%  '

$ws.Range('D15').Value = '3.530.97'
$ws.Range('E15').Value = '  -7.55%  '

$ws.Range('D16').Value = '64.357.34'
$ws.Range('E16').Value = '  -4.78%  '

$ws.Range('D17').Value = '0.111'
$ws.Range('E17').Value = '  -2.83%  '

$ws.Range('D18').Value = '3.051.06'
$ws.Range('E18').Value = '  -7.13%  '

$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '6.71'
$ws.Range('E19').Value = '  -8.31%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '489.07'
$ws.Range('E20').Value = '  -9.65%  '

$ws.Range('D21').Value = '13.67'
$ws.Range('E21').Value = '  -10.35%  '

$ws.Range('D22').Value = '0.687'
$ws.Range('E22').Value = '  -10.71%  '

$ws.Range('D23').Value = '7.12'
$ws.Range('E23').Value = '  -10.62%  '

$ws.Range('D24').Value = '12.65'
$ws.Range('E24').Value = '  -7.76%  '

$ws.Range('D25').Value = '78.81'
$ws.Range('E25').Value = '  -8.76%  '

$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.20%  '

$ws.Range('D27').Value = '2.74'
$ws.Range('E27').Value = '  -16.28%  '

$ws.Range('D28').Value = '2.11'
$ws.Range('E28').Value = '  -2.99%  '

$ws.Range('D29').Value = '7.77'
$ws.Range('E29').Value = '  -5.08%  '

$ws.Range('E30').Value = '  +0.14%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '25.93'
$ws.Range('E31').Value = '  -12.61%  '

$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').Value = '2.59'
$ws.Range('E32').Value = '  -3.89%  '

$ws.Range('D33').Value = '1.11'
$ws.Range('E33').Value = '  -3.48%  '

$ws.Range('D34').Value = '523.14'
$ws.Range('E34').Value = '  -4.65%  '

$ws.Range('D35').Value = '5.47'
$ws.Range('E35').Value = '  -6.00%  '

$ws.Range('D36').Value = '5.96'
$ws.Range('E36').Value = '  -10.89%  '

$ws.Range('D37').Value = '52.63'
$ws.Range('E37').Value = '  -2.03%  '

$ws.Range('D38').Value = '0.0405'
$ws.Range('E38').Value = '  -13.46%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.0799'
$ws.Range('E39').Value = '  -8.15%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.121'
$ws.Range('E40').Value = '  -6.20%  '

$ws.Range('D41').Value = '8.33'
$ws.Range('E41').Value = '  -9.34%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '2.73'
$ws.Range('E42').Value = '  -1.56%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.880.99'
$ws.Range('E43').Value = '  -2.62%  '

$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.09%  '

$ws.Range('D45').Value = '0.244'
$ws.Range('E45').Value = '  -9.91%  '

$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').Value = '2.08'
$ws.Range('E46').Value = '  -6.13%  '

$ws.Range('D47').Value = '0.0₃0534'
$ws.Range('E47').Value = '  -10.91%  '

$ws.Range('D48').Value = '24.95'
$ws.Range('E48').Value = '  -7.62%  '

$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '120.99'
$ws.Range('E49').Value = '  -4.50%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.108'
$ws.Range('E50').Value = '  -5.45%  '

$ws.Range('D51').Value = '2.04'
$ws.Range('E51').Value = '  -14.21%  '
